# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (stock) sheet gains three new trailing columns (H, I, J):
#   date, legislator_name, legislator_id
# with the single data row filled in from the source filename
# (林國正_2012-04-06, legislator id 1742).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Copy the existing header/data-row formatting into the new columns
# before writing any values, so H1:J1 pick up the bold/bordered header
# look (same as B1:G1) and H2:J2 pick up the plain data-row look (same
# as B2:G2), instead of ending up with Excel's bare default formatting.
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("H2:J2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New header row (row 1)
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# New data row (row 2)
# H2 is a literal date string ("2012-04-06"), not a real Excel date, so
# force text formatting before assigning it - otherwise Excel silently
# reinterprets it as a date serial number. Re-apply the plain data-row
# format afterwards so the cell matches its neighbours again.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "2012-04-06"
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I2").Value = "林國正"
$ws.Range("J2").Value = 1742
